$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 15
$ws.Range("H15").Value = 1186.63
$ws.Range("I15").Value = 1186.63
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 3559.89
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -3390.89
# ALC row 62
$ws.Range("H62").Value = 111111110
$ws.Range("I62").Value = 111111110
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 111111110
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -111110486
$ws.Range("N62").ClearContents()
# ALC row 65
$ws.Range("H65").Value = 111111110
$ws.Range("I65").Value = 111111110
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 555555550
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -555552430
$ws.Range("N65").ClearContents()
# ALC row 129
$ws.Range("H129").Value = 776.4474
$ws.Range("I129").Value = 375.2857
$ws.Range("J129").Value = 867.0323
$ws.Range("K129").Value = 1125.8571
$ws.Range("L129").Value = 2601.0969
$ws.Range("M129").Value = 3874.1429
$ws.Range("N129").Value = -12601.0969
# ALC row 133
$ws.Range("H133").Value = 34635
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 34635
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 34635
$ws.Range("N133").Value = -44755
# ALC row 137
$ws.Range("H137").Value = 1472.2808
$ws.Range("I137").Value = 1373.2122
$ws.Range("J137").Value = 1608.5
$ws.Range("K137").Value = 4119.6366
$ws.Range("L137").Value = 4825.5
$ws.Range("M137").Value = -1569.6366
$ws.Range("N137").Value = -9925.5

$ws = $wb.Worksheets.Item("ARM")
# ARM row 2
$ws.Range("H2").Value = 910.2308
$ws.Range("I2").Value = 728.1053000000001
$ws.Range("J2").Value = 1404.5714
$ws.Range("K2").Value = 728.1053000000001
$ws.Range("L2").Value = 1404.5714
$ws.Range("M2").Value = -615.1053000000001
$ws.Range("N2").Value = -1630.5714
# ARM row 5
$ws.Range("H5").Value = 196.23077
$ws.Range("I5").Value = 160.1
$ws.Range("J5").Value = 316.66666
$ws.Range("K5").Value = 160.1
$ws.Range("L5").Value = 316.66666
$ws.Range("M5").Value = -48.09999999999999
$ws.Range("N5").Value = -540.66666
# ARM row 32
$ws.Range("H32").Value = 4103.63
$ws.Range("I32").Value = 3776.9792
$ws.Range("J32").Value = 11943.25
$ws.Range("K32").Value = 3776.9792
$ws.Range("L32").Value = 11943.25
$ws.Range("M32").Value = -3489.9792
$ws.Range("N32").Value = -12517.25
# ARM row 74
$ws.Range("H74").Value = 3659.6
$ws.Range("I74").Value = 3260
$ws.Range("J74").Value = 4116.2856
$ws.Range("K74").Value = 3260
$ws.Range("L74").Value = 4116.2856
$ws.Range("M74").Value = -2386
$ws.Range("N74").Value = -5864.2856
# ARM row 77
$ws.Range("H77").Value = 3659.6
$ws.Range("I77").Value = 3260
$ws.Range("J77").Value = 4116.2856
$ws.Range("K77").Value = 16300
$ws.Range("L77").Value = 20581.428
$ws.Range("M77").Value = -11932
$ws.Range("N77").Value = -29317.428
# ARM row 116
$ws.Range("H116").Value = 910.2308
$ws.Range("I116").Value = 728.1053000000001
$ws.Range("J116").Value = 1404.5714
$ws.Range("K116").Value = 728.1053000000001
$ws.Range("L116").Value = 1404.5714
$ws.Range("M116").Value = 1565.8947
$ws.Range("N116").Value = -5992.5714
# ARM row 139
$ws.Range("H139").Value = 33289.375
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 33289.375
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 33289.375
$ws.Range("N139").Value = -43569.375

$ws = $wb.Worksheets.Item("BSM")
# BSM row 3
$ws.Range("H3").Value = 910.2308
$ws.Range("I3").Value = 728.1053000000001
$ws.Range("J3").Value = 1404.5714
$ws.Range("K3").Value = 728.1053000000001
$ws.Range("L3").Value = 1404.5714
$ws.Range("M3").Value = -614.1053000000001
$ws.Range("N3").Value = -1632.5714
# BSM row 4
$ws.Range("H4").Value = 196.23077
$ws.Range("I4").Value = 160.1
$ws.Range("J4").Value = 316.66666
$ws.Range("K4").Value = 160.1
$ws.Range("L4").Value = 316.66666
$ws.Range("M4").Value = -45.09999999999999
$ws.Range("N4").Value = -546.66666
# BSM row 94
$ws.Range("H94").Value = 7353282.5
$ws.Range("I94").Value = 8333632
$ws.Range("J94").Value = 659.75
$ws.Range("K94").Value = 8333632
$ws.Range("L94").Value = 659.75
$ws.Range("M94").Value = -8333181
$ws.Range("N94").Value = -1561.75
# BSM row 99
$ws.Range("H99").Value = 76924140
$ws.Range("I99").Value = 90910080
$ws.Range("J99").Value = 1490
$ws.Range("K99").Value = 90910080
$ws.Range("L99").Value = 1490
$ws.Range("M99").Value = -90908582
# BSM row 107
$ws.Range("H107").Value = 876.71875
$ws.Range("I107").Value = 737.15
$ws.Range("J107").Value = 1109.3334
$ws.Range("K107").Value = 737.15
$ws.Range("L107").Value = 1109.3334
$ws.Range("M107").Value = 1182.85

$ws = $wb.Worksheets.Item("CRP")
# CRP row 7
$ws.Range("H7").Value = 215.04347
$ws.Range("I7").Value = 86.90909000000001
$ws.Range("J7").Value = 332.5
$ws.Range("K7").Value = 86.90909000000001
$ws.Range("L7").Value = 332.5
$ws.Range("M7").Value = 26.09090999999999
$ws.Range("N7").Value = -558.5
# CRP row 16
$ws.Range("H16").Value = 90910660
$ws.Range("I16").Value = 111112740
$ws.Range("J16").Value = 1325
$ws.Range("K16").Value = 111112740
$ws.Range("L16").Value = 1325
$ws.Range("M16").Value = -111112453
$ws.Range("N16").Value = -1899
# CRP row 31
$ws.Range("H31").Value = 1270.7301
$ws.Range("I31").Value = 1236.7451
$ws.Range("J31").Value = 1415.1666
$ws.Range("K31").Value = 1236.7451
$ws.Range("L31").Value = 1415.1666
$ws.Range("M31").Value = -941.7451000000001
$ws.Range("N31").Value = -2005.1666
# CRP row 34
$ws.Range("H34").Value = 1270.7301
$ws.Range("I34").Value = 1236.7451
$ws.Range("J34").Value = 1415.1666
$ws.Range("K34").Value = 1236.7451
$ws.Range("L34").Value = 1415.1666
$ws.Range("M34").Value = -1034.7451
$ws.Range("N34").Value = -1819.1666
# CRP row 107
$ws.Range("H107").Value = 875.2353000000001
$ws.Range("I107").Value = 452.23077
$ws.Range("J107").Value = 2250
$ws.Range("K107").Value = 452.23077
$ws.Range("L107").Value = 2250
$ws.Range("M107").Value = 1467.76923
$ws.Range("N107").Value = -6090
# CRP row 113
$ws.Range("H113").Value = 90910660
$ws.Range("I113").Value = 111112740
$ws.Range("J113").Value = 1325
$ws.Range("K113").Value = 111112740
$ws.Range("L113").Value = 1325
$ws.Range("M113").Value = -111110570
$ws.Range("N113").Value = -5665

$ws = $wb.Worksheets.Item("CUL")
# CUL row 39
$ws.Range("H39").Value = 4038.3845
$ws.Range("I39").Value = 3800
$ws.Range("J39").Value = 4081.7273
$ws.Range("K39").Value = 11400
$ws.Range("L39").Value = 12245.1819
$ws.Range("M39").Value = -11106
$ws.Range("N39").Value = -12833.1819
# CUL row 98
$ws.Range("H98").Value = 1639.5
$ws.Range("I98").Value = 2271.4
$ws.Range("J98").Value = 586.3333
$ws.Range("K98").Value = 6814.200000000001
$ws.Range("L98").Value = 1758.9999
$ws.Range("M98").Value = -5316.200000000001
# CUL row 107
$ws.Range("H107").Value = 3734.8667
$ws.Range("I107").Value = 496
$ws.Range("J107").Value = 4382.64
$ws.Range("K107").Value = 1488
$ws.Range("L107").Value = 13147.92
$ws.Range("M107").Value = 432
$ws.Range("N107").Value = -16987.92

$ws = $wb.Worksheets.Item("GSM")
# GSM row 113
$ws.Range("H113").Value = 1887
$ws.Range("I113").Value = 1887
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1887
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 283
$ws.Range("N113").ClearContents()
# GSM row 132
$ws.Range("H132").Value = 4187.6113
$ws.Range("I132").Value = 4705.6665
$ws.Range("J132").Value = 3669.5557
$ws.Range("K132").Value = 14116.9995
$ws.Range("L132").Value = 11008.6671
$ws.Range("M132").Value = -11586.9995
$ws.Range("N132").Value = -16068.6671

$ws = $wb.Worksheets.Item("LTW")
# LTW row 20
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
# LTW row 132
$ws.Range("H132").Value = 3778.2856
$ws.Range("I132").Value = 4879.6
$ws.Range("J132").Value = 3166.4443
$ws.Range("K132").Value = 14638.8
$ws.Range("L132").Value = 9499.332900000001
$ws.Range("M132").Value = -12108.8
$ws.Range("N132").Value = -14559.3329

$ws = $wb.Worksheets.Item("WVR")
# WVR row 62
$ws.Range("H62").Value = 500000000
$ws.Range("I62").Value = 500000000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 500000000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -499999376
$ws.Range("N62").ClearContents()
# WVR row 65
$ws.Range("H65").Value = 500000000
$ws.Range("I65").Value = 500000000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 2500000000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -2499996880
$ws.Range("N65").ClearContents()
# WVR row 81
$ws.Range("H81").Value = 1527.2858
$ws.Range("I81").Value = 1198.5
$ws.Range("J81").Value = 3500
$ws.Range("K81").Value = 2397
$ws.Range("L81").Value = 7000
$ws.Range("M81").Value = -1336
$ws.Range("N81").Value = -9122
# WVR row 84
$ws.Range("H84").Value = 1527.2858
$ws.Range("I84").Value = 1198.5
$ws.Range("J84").Value = 3500
$ws.Range("K84").Value = 11985
$ws.Range("L84").Value = 35000
$ws.Range("M84").Value = -6681
$ws.Range("N84").Value = -45608
# WVR row 136
$ws.Range("H136").Value = 1538.8
$ws.Range("I136").Value = 1379.8889
$ws.Range("J136").Value = 1668.8182
$ws.Range("K136").Value = 4139.6667
$ws.Range("L136").Value = 5006.4546
$ws.Range("M136").Value = -1589.6667
